# Updated cryptos list - apply price / 1h volume changes from commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.964.87"
$ws.Range("E2").Value = "'  +0.06%  "
$ws.Range("D3").Value = "'3.384.33"
$ws.Range("E3").Value = "'  -0.97%  "
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("D5").Value = "'570.93"
$ws.Range("E5").Value = "'  -0.19%  "
$ws.Range("D6").Value = "'141.65"
$ws.Range("E6").Value = "'  -0.36%  "
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E8").Value = "'  -0.53%  "
$ws.Range("E9").Value = "'  +0.64%  "
$ws.Range("E10").Value = "'  -1.63%  "
$ws.Range("E11").Value = "'  +1.26%  "
$ws.Range("D12").Value = "'3.964.24"
$ws.Range("E12").Value = "'  -1.00%  "
$ws.Range("E13").Value = "'  +2.33%  "
$ws.Range("D14").Value = "'27.62"
$ws.Range("E14").Value = "'  -2.31%  "
$ws.Range("E15").Value = "'  -0.95%  "
$ws.Range("D16").Value = "'3.374.61"
$ws.Range("E16").Value = "'  -1.59%  "
$ws.Range("D17").Value = "'61.061.26"
$ws.Range("E17").Value = "'  -0.10%  "
$ws.Range("D18").Value = "'6.09"
$ws.Range("E18").Value = "'  -4.21%  "
$ws.Range("D19").Value = "'13.70"
$ws.Range("E19").Value = "'  -4.81%  "
$ws.Range("E20").Value = "'  -4.90%  "
$ws.Range("D21").Value = "'381.16"
$ws.Range("E21").Value = "'  -4.03%  "
$ws.Range("D22").Value = "'74.83"
$ws.Range("E22").Value = "'  +2.39%  "
$ws.Range("D23").Value = "'0.553"
$ws.Range("E23").Value = "'  -1.90%  "
$ws.Range("E24").Value = "'  +0.44%  "
$ws.Range("D25").Value = "'0.0000117"
$ws.Range("E25").Value = "'  -4.03%  "
$ws.Range("D26").Value = "'3.525.38"
$ws.Range("E26").Value = "'  -1.42%  "
$ws.Range("E27").Value = "'  +1.56%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "'  +0.03%  "
$ws.Range("E29").Value = "'  -1.99%  "
$ws.Range("E30").Value = "'  -0.75%  "
$ws.Range("D31").Value = "'7.95"
$ws.Range("E31").Value = "'  -2.32%  "
$ws.Range("D32").Value = "'1.40"
$ws.Range("E32").Value = "'  -2.23%  "
$ws.Range("E33").Value = "'  -0.05%  "
$ws.Range("D34").Value = "'23.35"
$ws.Range("E34").Value = "'  -2.10%  "
$ws.Range("E35").Value = "'  -0.64%  "
$ws.Range("D36").Value = "'166.17"
$ws.Range("E36").Value = "'  -0.58%  "
$ws.Range("D37").Value = "'3.416.99"
$ws.Range("E37").Value = "'  -0.86%  "
$ws.Range("E38").Value = "'  -2.25%  "
$ws.Range("E39").Value = "'  -5.08%  "
$ws.Range("D40").Value = "'0.0767"
$ws.Range("E40").Value = "'  -2.27%  "
$ws.Range("D41").Value = "'27.29"
$ws.Range("E41").Value = "'  +0.01%  "
$ws.Range("E42").Value = "'  -0.14%  "
$ws.Range("E43").Value = "'  -2.40%  "
$ws.Range("D44").Value = "'41.81"
$ws.Range("E44").Value = "'  -0.74%  "
$ws.Range("E45").Value = "'  -2.36%  "
$ws.Range("E46").Value = "'  -2.94%  "
$ws.Range("E47").Value = "'  -0.93%  "
$ws.Range("D48").Value = "'2.454.34"
$ws.Range("E48").Value = "'  -6.08%  "
$ws.Range("D49").Value = "'22.96"
$ws.Range("E49").Value = "'  +0.27%  "
$ws.Range("D50").Value = "'6.74"
$ws.Range("E50").Value = "'  -3.11%  "
$ws.Range("D51").Value = "'0.0266"
$ws.Range("E51").Value = "'  +2.08%  "
